$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-39 changed from serial date 45182 to 45184
for ($row = 2; $row -le 39; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
